$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, centered) from E1 onto the new header cells F1:H1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New header text for F1:H1
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# New boolean (FALSE) values for the data rows 2-5, columns F:H
foreach ($r in 2..5) {
    $ws.Range("F$r").Value = $false
    $ws.Range("G$r").Value = $false
    $ws.Range("H$r").Value = $false
}
